$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.341.57'
$ws.Range("E2").Value = '  +1.71%  '

$ws.Range("D3").Value = '2.507.62'
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.84%  '

$ws.Range("E7").Value = '  +1.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.541'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("E10").Value = '  +1.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.32'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.97%  '

$ws.Range("E12").Value = '  +1.15%  '

$ws.Range("E13").Value = '  -0.03%  '

$ws.Range("E14").Value = '  +0.10%  '

$ws.Range("D15").Value = '2.899.46'
$ws.Range("E15").Value = '  +0.62%  '

$ws.Range("D16").Value = '2.505.40'
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("E17").Value = '  -0.44%  '

$ws.Range("D18").Value = '48.185.47'
$ws.Range("E18").Value = '  +1.63%  '

$ws.Range("E19").Value = '  -2.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.57%  '

$ws.Range("E21").Value = '  +1.29%  '

$ws.Range("D22").Value = '0.0₃0947'
$ws.Range("E22").Value = '  +0.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '281.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +13.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.28%  '

$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.83%  '

$ws.Range("E30").Value = '  +1.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.38'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.85%  '

$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("E36").Value = '  -0.45%  '

$ws.Range("E37").Value = '  -0.30%  '

$ws.Range("E38").Value = '  -1.85%  '

$ws.Range("E39").Value = '  -0.27%  '

$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '121.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.17%  '

$ws.Range("E44").Value = '  +2.23%  '

$ws.Range("D45").Value = '2.012.59'
$ws.Range("E45").Value = '  +0.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.42%  '

$ws.Range("E47").Value = '  +3.65%  '

$ws.Range("E48").Value = '  -2.77%  '

$ws.Range("E49").Value = '  -0.71%  '

$ws.Range("E50").Value = '  -0.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.95%  '
